$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update team header to include Frank ---
$ws.Range("C1").Value = "TEAM:Donnovan,Daniel,Frank"

# --- 2. Fix row 7 (Consultar Producto) function name ---
$ws.Range("C10").Value = "f_consultar_producto"

# --- 3. Build new "persona" CRUD rows (15-19) by copying the formatting
#        of the existing "producto" row 10 (the "plain"/non-last-row style) ---
for ($r = 15; $r -le 19; $r++) {
    $ws.Range("A10:G10").Copy($ws.Range("A" + $r + ":G" + $r))
}

# Row 19 is the last row of the block: column F there uses the special
# last-row date-format style, so copy that single cell's formatting from F14.
$ws.Range("F14").Copy($ws.Range("F19"))

# --- 4. Recolor the new block's fill to the Accent6 (green) theme color ---
$fillRange = $ws.Range("A15:G19")
$fillRange.Interior.ThemeColor = 10

# --- 5. Setting the theme color can corrupt the builtin date number format
#        on F19 (turns "m/d/yyyy" into a custom format) - restore it. ---
$ws.Range("F19").NumberFormat = "mm-dd-yy"

# --- 6. Set the explicit row heights for the new rows ---
$ws.Rows.Item(15).RowHeight = 21
$ws.Rows.Item(16).RowHeight = 21
$ws.Rows.Item(17).RowHeight = 21
$ws.Rows.Item(18).RowHeight = 21
$ws.Rows.Item(19).RowHeight = 21

# --- 7. Fill in the values for the new "persona" CRUD rows ---
$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "Consultar Producto"
$ws.Range("C15").Value = "f_consultar_persona"
$ws.Range("D15").Value = "frank"
$ws.Range("E15").Value = "develop"
$ws.Range("G15").Value = "frank"

$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "Agregar Producto"
$ws.Range("C16").Value = "f_agregar_persona"
$ws.Range("D16").Value = "frank"
$ws.Range("E16").Value = "develop"
$ws.Range("G16").Value = "frank"

$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "Eliminar Producto"
$ws.Range("C17").Value = "f_eliminar_persona"
$ws.Range("D17").Value = "frank"
$ws.Range("E17").Value = "develop"
$ws.Range("G17").Value = "frank"

$ws.Range("A18").Value = 15
$ws.Range("B18").Value = "Actualizar Producto"
$ws.Range("C18").Value = "f_actualizar_persona"
$ws.Range("D18").Value = "frank"
$ws.Range("E18").Value = "develop"
$ws.Range("G18").Value = "frank"

$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "Crear Menu de Opciones"
$ws.Range("C19").Value = "f_menu_persona"
$ws.Range("D19").Value = "frank"
$ws.Range("E19").Value = "develop"
$ws.Range("G19").Value = "frank"

# --- 8. Update the active selection to match the final workbook state ---
$ws.Range("E21").Select()
